$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three formula-text strings in column E (rows 2-4)
$ws.Range("E2").Value = "3*10*100*100=300000"
$ws.Range("E3").Value = "125*2*4*5*100=500000"
$ws.Range("E4").Value = "40*5*10*100=200000"

# Update the active cell selection to F8
$ws.Range("F8").Select()
